$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings (e.g. "1.000",
# "30.323.24") are preserved verbatim as text rather than coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.323.24'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '1.880.61'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '237.45'
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.4824'
$ws.Range("E7").Value = '  -2.29%  '
$ws.Range("D8").Value = '0.2891'
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").Value = '0.06587'
$ws.Range("E9").Value = '  -2.71%  '
$ws.Range("D10").Value = '1.879.76'
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("D11").Value = '16.93'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '0.07392'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '5.190'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = '88.03'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").Value = '0.6599'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '30.288.96'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").Value = '13.62'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").Value = '0.9998'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '0.000007728'
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("D20").Value = '5.473'
$ws.Range("E20").Value = '  +2.18%  '
$ws.Range("D21").Value = '2.137.74'
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '195.95'
$ws.Range("E23").Value = '  -2.40%  '
$ws.Range("D24").Value = '6.159'
$ws.Range("E24").Value = '  -2.38%  '
$ws.Range("D25").Value = '9.436'
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("D26").Value = '163.37'
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("D27").Value = '18.25'
$ws.Range("E27").Value = '  -3.35%  '
$ws.Range("D28").Value = '1.928'
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("D29").Value = '1.441'
$ws.Range("E29").Value = '  -2.67%  '
$ws.Range("D30").Value = '4.278'
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").Value = '0.09146'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").Value = '4.052'
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("D33").Value = '0.05047'
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("D34").Value = '0.7416'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Value = '1.141'
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("D36").Value = '2.707'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").Value = '0.01841'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '2.632'
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("D39").Value = '0.9161'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '2.072'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").Value = '106.36'
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("D42").Value = '0.4321'
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("D43").Value = '5.868'
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '7.638'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '0.1349'
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("D47").Value = '1.570'
$ws.Range("E47").Value = '  +9.18%  '
$ws.Range("D48").Value = '65.11'
$ws.Range("E48").Value = '  -11.50%  '
$ws.Range("D49").Value = '8.890'
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("D50").Value = '34.18'
$ws.Range("E50").Value = '  -3.83%  '
$ws.Range("D51").Value = '0.05722'
$ws.Range("E51").Value = '  -2.69%  '
